$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a text value to a cell while forcing text storage
# (leading apostrophe = Excel quote-prefix) without leaving a residual
# "Text" number-format style behind on the cell.
function Set-TextCell($cellRef, $text) {
    $ws.Range($cellRef).Value = "`'$text"
    $ws.Range($cellRef).Style = "Normal"
}

Set-TextCell "D2" "26.493.56"
Set-TextCell "E2" "  -7.72%  "

Set-TextCell "D3" "1.678.80"
Set-TextCell "E3" "  -6.91%  "

Set-TextCell "D4" "1.005"
Set-TextCell "E4" "  +0.27%  "

Set-TextCell "D5" "216.30"
Set-TextCell "E5" "  -6.59%  "

Set-TextCell "E6" "  +0.13%  "

Set-TextCell "D7" "0.4960"
Set-TextCell "E7" "  -16.62%  "

Set-TextCell "D8" "0.2596"
Set-TextCell "E8" "  -6.72%  "

Set-TextCell "D9" "21.75"
Set-TextCell "E9" "  -6.73%  "

Set-TextCell "D10" "0.06177"
Set-TextCell "E10" "  -9.66%  "

Set-TextCell "D11" "0.07275"
Set-TextCell "E11" "  -3.47%  "

Set-TextCell "D12" "1.624.50"
Set-TextCell "E12" "  -7.34%  "

Set-TextCell "D13" "4.425"
Set-TextCell "E13" "  -7.31%  "

Set-TextCell "D14" "0.5733"
Set-TextCell "E14" "  -8.13%  "

Set-TextCell "D15" "1.906.03"
Set-TextCell "E15" "  -6.96%  "

Set-TextCell "D16" "0.000008155"
Set-TextCell "E16" "  -12.70%  "

Set-TextCell "D17" "64.24"
Set-TextCell "E17" "  -14.98%  "

Set-TextCell "D18" "26.497.47"

Set-TextCell "D19" "4.974"
Set-TextCell "E19" "  -9.37%  "

Set-TextCell "E20" "  +0.16%  "

Set-TextCell "D21" "10.74"
Set-TextCell "E21" "  -6.23%  "

Set-TextCell "D22" "183.69"
Set-TextCell "E22" "  -12.55%  "

Set-TextCell "D23" "6.166"
Set-TextCell "E23" "  -10.12%  "

Set-TextCell "D24" "1.006"

Set-TextCell "D25" "144.31"
Set-TextCell "E25" "  -6.48%  "

Set-TextCell "D26" "7.443"
Set-TextCell "E26" "  -5.40%  "

Set-TextCell "D27" "0.1127"
Set-TextCell "E27" "  -11.42%  "

Set-TextCell "D28" "15.38"
Set-TextCell "E28" "  -6.08%  "

Set-TextCell "D29" "1.300"
Set-TextCell "E29" "  -9.03%  "

Set-TextCell "D30" "0.05673"
Set-TextCell "E30" "  -8.67%  "

Set-TextCell "D31" "1.317"
Set-TextCell "E31" "  -7.27%  "

Set-TextCell "D32" "3.466"
Set-TextCell "E32" "  -8.32%  "

Set-TextCell "D33" "3.455"
Set-TextCell "E33" "  -7.87%  "

Set-TextCell "D34" "1.631"
Set-TextCell "E34" "  -5.26%  "

Set-TextCell "D35" "1.003"
Set-TextCell "E35" "  -5.75%  "

Set-TextCell "D36" "2.367"
Set-TextCell "E36" "  -5.05%  "

Set-TextCell "D37" "0.5886"
Set-TextCell "E37" "  -7.95%  "

Set-TextCell "D38" "2.630"
Set-TextCell "E38" "  -3.04%  "

Set-TextCell "D39" "0.01584"
Set-TextCell "E39" "  -7.67%  "

Set-TextCell "D40" "1.067.56"
Set-TextCell "E40" "  -5.73%  "

Set-TextCell "E41" "  -9.20%  "

Set-TextCell "E42" "  -3.04%  "

Set-TextCell "D43" "1.001"
Set-TextCell "E43" "  -0.32%  "

Set-TextCell "D44" "98.11"
Set-TextCell "E44" "  -2.71%  "

Set-TextCell "D45" "1.835.68"
Set-TextCell "E45" "  -6.35%  "

Set-TextCell "D46" "56.00"
Set-TextCell "E46" "  -7.52%  "

# Row 47: coin swap
$ws.Range("B47").Value = "Frax"
$ws.Range("C47").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
Set-TextCell "D47" "1.004"
Set-TextCell "E47" "  -0.33%  "

# Row 48: coin swap
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextCell "D48" "0.00000000104"
Set-TextCell "E48" "  -7.06%  "

Set-TextCell "D49" "7.984"
Set-TextCell "E49" "  -4.60%  "

Set-TextCell "D50" "0.4310"
Set-TextCell "E50" "  -3.87%  "

Set-TextCell "D51" "0.05185"
Set-TextCell "E51" "  -5.32%  "
